$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.811038136482239
$ws.Range("B1").Value = 4.736396312713623
$ws.Range("C1").Value = 3.769814968109131
$ws.Range("D1").Value = 0.9020178318023682
$ws.Range("E1").Value = 0.4737375676631927
